$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header F1 text
$ws.Range("F1").Value = "Firma.1"

# Add new attendance row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "ALBARRAN  JIMENEZ OSCAR EDUARDO"
$ws.Range("C3").Value = "23:42:00"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "23:42:00"
$ws.Range("F3").ClearFormats()
